$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 186. This shifts the existing rows
# 186-191 down to 188-193, preserving their data/styles.
$insertRange = $ws.Range("A186:R187")
$insertRange.EntireRow.Insert()

# Populate the two newly inserted rows (186 and 187) with the new week's data.
# Row 186 - "Primera" quality
$ws.Cells.Item(186, 1).Value = 11
$ws.Cells.Item(186, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(186, 3).Value = "Bíobío"
$ws.Cells.Item(186, 4).Value = 44615
$ws.Cells.Item(186, 4).NumberFormat = $ws.Cells.Item(190, 4).NumberFormat
$ws.Cells.Item(186, 5).Value = 8
$ws.Cells.Item(186, 6).Value = 100114013
$ws.Cells.Item(186, 7).Value = "Zanahoria"
$ws.Cells.Item(186, 8).Value = "Sin especificar"
$ws.Cells.Item(186, 9).Value = "Primera"
$ws.Cells.Item(186, 10).Value = 1000
$ws.Cells.Item(186, 11).Value = 8000
$ws.Cells.Item(186, 12).Value = 8500
$ws.Cells.Item(186, 13).Value = 8250
$ws.Cells.Item(186, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(186, 15).Value = "Región de Ñuble"
$ws.Cells.Item(186, 16).Value = 412
$ws.Cells.Item(186, 17).Value = 20
$ws.Cells.Item(186, 18).Value = "Hortaliza"

# Row 187 - "Segunda" quality
$ws.Cells.Item(187, 1).Value = 11
$ws.Cells.Item(187, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(187, 3).Value = "Bíobío"
$ws.Cells.Item(187, 4).Value = 44615
$ws.Cells.Item(187, 4).NumberFormat = $ws.Cells.Item(190, 4).NumberFormat
$ws.Cells.Item(187, 5).Value = 8
$ws.Cells.Item(187, 6).Value = 100114013
$ws.Cells.Item(187, 7).Value = "Zanahoria"
$ws.Cells.Item(187, 8).Value = "Sin especificar"
$ws.Cells.Item(187, 9).Value = "Segunda"
$ws.Cells.Item(187, 10).Value = 500
$ws.Cells.Item(187, 11).Value = 7000
$ws.Cells.Item(187, 12).Value = 7000
$ws.Cells.Item(187, 13).Value = 7000
$ws.Cells.Item(187, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(187, 15).Value = "Región de Ñuble"
$ws.Cells.Item(187, 16).Value = 350
$ws.Cells.Item(187, 17).Value = 20
$ws.Cells.Item(187, 18).Value = "Hortaliza"
